# "Generate Report for Archive"
# The file with uuid 4ea1af5a-5d3c-42fa-8e63-933cff256c94.md has moved on from
# "Ready for handoff" to "In Translation". Update its Status on every sheet
# that tracks it: the Overview roll-up sheet (one column per locale) and each
# locale-specific handoff sheet (a single Status column).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 4 is the 4ea1af5a...md file; B (zh-cn) & C (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

# --- Per-locale sheets: column B is "Status"; row 4 is the same file ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $locSheet = $wb.Worksheets.Item($sheetName)
    $locSheet.Range("B4").Value = "In Translation"
}
